$d = $word.ActiveDocument

function New-WordXmlPackage([string]$InnerBodyXml) {
    $ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $ns + '><w:body>' + $InnerBodyXml + '</w:body></w:document></pkg:xmlData>' +
        '</pkg:part></pkg:package>'
}

# --- 1. Paragraph "Todays date is 26 november" -> "Todays date is 1 november",
#        re-typed as several spell-checked runs (matches the proofErr markup
#        Word leaves behind for words it doesn't recognise). -----------------
$needle = $d.Content
$foundIt = $needle.Find.Execute("Todays date is 26 november")
if (-not $foundIt) {
    throw "Could not locate the target phrase 'Todays date is 26 november'"
}
$needleStart = $needle.Start

# resolve the paragraph that actually contains the found text (locating by
# position rather than a hard-coded index keeps this robust to reordering)
$paraIndex = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pr = $d.Paragraphs.Item($i).Range
    if ($needleStart -ge $pr.Start -and $needleStart -lt $pr.End) {
        $paraIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($paraIndex).Range
$origStart = $target.Start
$origEnd = $target.End
$origLen = $origEnd - $origStart

$newRunsXml =
    '<w:p>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Todays</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> date is</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> 1</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>november</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

$insertionPoint = $d.Range($origStart, $origStart)
$insertionPoint.InsertXML((New-WordXmlPackage $newRunsXml))

# the freshly-typed text now sits right before the untouched original text;
# trim the stale original run(s) off the tail of the (same) paragraph.
$grownEnd = $d.Paragraphs.Item($paraIndex).Range.End
$staleRange = $d.Range($grownEnd - $origLen, $grownEnd)
$staleRange.Delete()

# --- 2. Two brand-new paragraphs appended at the end of the document. ------
$docEnd = $d.Range($d.Content.End, $d.Content.End)

$timeParaXml =
    '<w:p>' +
        '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Time now is 19:22</w:t></w:r>' +
    '</w:p>'
$docEnd.InsertXML((New-WordXmlPackage $timeParaXml))

$docEnd = $d.Range($d.Content.End, $d.Content.End)
$githubParaXml =
    '<w:p>' +
        '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">We are learning </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>github</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
$docEnd.InsertXML((New-WordXmlPackage $githubParaXml))

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host ("Para {0}: [{1}]" -f $i, $d.Paragraphs.Item($i).Range.Text)
}
